$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra "full block" character-design example block (columns BN:BT)
# which lived alongside the 8 character blocks in rows 39-47, plus the spacer
# cell BE40:BE47 that preceded it. This also drops the now-unused shared
# string describing that example block.
$ws.Range("BN39:BT47").Clear()
$ws.Range("BE40:BE47").Clear()

# The removed block was a solid "full block" character (all bits sets to 1).
# Re-use that design for the previously-blank character occupying the
# diagonal set of 5-column blocks in rows 40-47 (one block per row), instead
# of keeping it as a separate, duplicate block.
$fmtSrc = $ws.Range("BG47")

$targets = @(
    @{ Row = 40; Cols = "BG:BK" },
    @{ Row = 41; Cols = "AY:BC" },
    @{ Row = 42; Cols = "AQ:AU" },
    @{ Row = 43; Cols = "AI:AM" },
    @{ Row = 44; Cols = "AA:AE" },
    @{ Row = 45; Cols = "S:W" },
    @{ Row = 46; Cols = "K:O" },
    @{ Row = 47; Cols = "C:G" }
)

foreach ($t in $targets) {
    $parts = $t.Cols -split ":"
    $left = $parts[0]
    $right = $parts[1]
    $rng = $ws.Range("$left$($t.Row):$right$($t.Row)")
    $fmtSrc.Copy()
    $rng.PasteSpecial(-4122)
    $rng.Value = 1
}

$excel.CutCopyMode = 0

# Match the resulting selection left after the edit.
$ws.Range("A40").Select()
